# Fix missing response validation:
#  - customer-accounts: renumber account-alias values for usr03's saving
#    account and all of usr04 / usr05's accounts (the old "checking6" /
#    "saving6" rows go away, shifting the remaining aliases down by one)
#  - verifyTransferPositiveCases: row 6 ("rounding is NOT a ceil rounding")
#    now transfers from usr01's own default/checking1 accounts instead of
#    the (removed) usr03 saving3 account
#  - users: usr101 is replaced by usr04
#  - selection / active-sheet bookkeeping to match the saved UI state

$wb = $excel.ActiveWorkbook

$wsAccounts = $wb.Worksheets.Item("customer-accounts")
$wsPos      = $wb.Worksheets.Item("verifyTransferPositiveCases")
$wsNeg      = $wb.Worksheets.Item("verifyTransferNegativeCases")
$wsUsers    = $wb.Worksheets.Item("users")

# --- customer-accounts: shift account-alias names down for usr03(saving)/usr04/usr05 ---
$wsAccounts.Range("B7").Value  = "saving3"
$wsAccounts.Range("B8").Value  = "checking4"
$wsAccounts.Range("B9").Value  = "saving4"
$wsAccounts.Range("B10").Value = "checking5"
$wsAccounts.Range("B11").Value = "saving5"

# --- verifyTransferPositiveCases: row 6 now uses usr01 / usr01.default / checking1 ---
$wsPos.Range("C6").Value = "usr01"
$wsPos.Range("D6").Value = "usr01.default"
$wsPos.Range("E6").Value = "checking1"

# --- users: usr101 -> usr04 ---
$wsUsers.Range("B5").Value = "usr04"

# --- view / selection state ---
$wsUsers.Range("B6").Select()

$wsNeg.Range("E2").Select()

$wsPos.Range("E2").Select()

$wsAccounts.Activate()
$excel.ActiveWindow.Zoom = 157
$wsAccounts.Range("B12").Select()

Write-Output "edit applied"
